$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.896.34'
$ws.Range("E2").Value = '  +0.24%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.889.67'
$ws.Range("E3").Value = '  +0.01%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7690'
$ws.Range("E5").Value = '  -1.01%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.80'
$ws.Range("E6").Value = '  -0.44%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3134'
$ws.Range("E8").Value = '  -0.14%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.69'
$ws.Range("E9").Value = '  +1.65%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07141'
$ws.Range("E10").Value = '  -2.62%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08535'
$ws.Range("E11").Value = '  +4.82%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7640'
$ws.Range("E12").Value = '  -0.14%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.371'
$ws.Range("E13").Value = '  -1.48%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.875.04'
$ws.Range("E14").Value = '  -1.48%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.88'
$ws.Range("E15").Value = '  +0.94%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.172'
$ws.Range("E16").Value = '  -0.54%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.863.82'
$ws.Range("E17").Value = '  +0.12%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.78'
$ws.Range("E18").Value = '  -0.95%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.54'
$ws.Range("E19").Value = '  -0.35%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007809'
$ws.Range("E20").Value = '  -0.44%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9990'
$ws.Range("E21").Value = '  -0.10%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.017'
$ws.Range("E22").Value = '  -1.56%  '

# Row 23
$ws.Range("E23").Value = '  -0.03%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1624'
$ws.Range("E24").Value = '  +3.45%  '

# Row 25
$ws.Range("E25").Value = '  +0.20%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.15'
$ws.Range("E26").Value = '  +1.38%  '

# Row 27
$ws.Range("E27").Value = '  +0.25%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.039'
$ws.Range("E28").Value = '  +0.20%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.502'
$ws.Range("E29").Value = '  +3.62%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.541'

# Row 31
$ws.Range("E31").Value = '  +0.41%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.115'
$ws.Range("E32").Value = '  +1.08%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05452'
$ws.Range("E33").Value = '  -2.16%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.242'
$ws.Range("E34").Value = '  -0.30%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7456'
$ws.Range("E35").Value = '  -1.11%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.001'
$ws.Range("E36").Value = '  +0.41%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.696'
$ws.Range("E37").Value = '  +2.26%  '

# Row 38
$ws.Range("E38").Value = '  +0.96%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.783'
$ws.Range("E39").Value = '  +0.31%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4472'
$ws.Range("E40").Value = '  +0.68%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.102.47'
$ws.Range("E41").Value = '  -3.57%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.090'
$ws.Range("E42").Value = '  +2.30%  '

# Row 43
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '73.14'
$ws.Range("E43").Value = '  -0.90%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8536'
$ws.Range("E44").Value = '  +0.26%  '

# Row 45
$ws.Range("E45").Value = '  -0.01%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.97'
$ws.Range("E46").Value = '  +1.07%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.870'
$ws.Range("E47").Value = '  -1.46%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.672'
$ws.Range("E48").Value = '  +2.48%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.058'
$ws.Range("E49").Value = '  -1.31%  '

# Row 50
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.004.75'
$ws.Range("E50").Value = '  -1.26%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06085'
$ws.Range("E51").Value = '  +0.61%  '
